$wb = $excel.ActiveWorkbook

# --- Sheet "Reguły": update rule description text for rows whose
#     object-id lists inside the brackets were reordered ---
$wsRules = $wb.Worksheets.Item("Reguły")
$wsRules.Range("B3").Value = '(age >=  39.0) & (infertility >=  3.0) => (class <= 2) [''a11'', ''a23'', ''a6'', ''a33'', ''a5'', ''a51'']'
$wsRules.Range("B5").Value = '(develop_quality <=  4.0) => (class <= 3) [''a26'', ''a9'', ''a33'', ''a38'', ''a34'', ''a14'', ''a51'', ''a29'', ''a49'', ''a12'', ''a30'', ''a17'', ''a24'', ''a23'', ''a25'', ''a43'', ''a5'', ''a21'', ''a4'', ''a44'', ''a47'', ''a16'', ''a50'', ''a7'', ''a11'', ''a45'', ''a41'', ''a1'', ''a3'', ''a13'', ''a19'', ''a22'']'
$wsRules.Range("B6").Value = '(infertility >=  4.0) => (class <= 3) [''a30'', ''a17'', ''a42'', ''a6'', ''a41'', ''a27'', ''a32'', ''a36'', ''a16'', ''a51'', ''a31'', ''a21'', ''a4'', ''a29'']'
$wsRules.Range("B7").Value = '(oocytes >=  4.0) => (class <= 3) [''a26'', ''a9'', ''a33'', ''a38'', ''a34'', ''a51'', ''a49'', ''a12'', ''a30'', ''a32'', ''a21'', ''a44'', ''a8'', ''a6'', ''a37'', ''a36'', ''a16'', ''a7'', ''a11'', ''a41'', ''a1'', ''a3'', ''a13'', ''a22'']'
$wsRules.Range("B8").Value = '(age >=  35.0) => (class <= 3) [''a9'', ''a33'', ''a34'', ''a51'', ''a29'', ''a49'', ''a12'', ''a17'', ''a23'', ''a5'', ''a48'', ''a6'', ''a36'', ''a11'', ''a39'', ''a41'', ''a1'', ''a27'', ''a3'', ''a13'', ''a22'']'
$wsRules.Range("B9").Value = '(sperm <=  2.0) => (class <= 3) [''a5'', ''a42'', ''a2'', ''a7'']'
$wsRules.Range("B10").Value = '(age <=  34.0) & (morpho_quality >=  10.0) => (class >= 3) [''a26'', ''a18'', ''a38'', ''a14'', ''a40'', ''a24'', ''a25'', ''a32'', ''a43'', ''a21'', ''a44'', ''a2'', ''a10'', ''a8'', ''a20'', ''a47'', ''a37'', ''a16'', ''a31'', ''a35'', ''a50'', ''a46'', ''a42'', ''a45'', ''a28'', ''a15'', ''a19'']'
$wsRules.Range("B11").Value = '(morpho_quality >=  15.0) & (oocytes <=  5.0) => (class >= 3) [''a18'', ''a38'', ''a34'', ''a14'', ''a29'', ''a2'', ''a10'', ''a8'', ''a20'', ''a47'', ''a48'', ''a37'', ''a16'', ''a31'', ''a50'', ''a39'', ''a46'', ''a42'', ''a27'', ''a13'', ''a15'']'
$wsRules.Range("B12").Value = '(age <=  36.0) & (infertility <=  3.0) & (oocytes <=  5.0) => (class >= 3) [''a26'', ''a9'', ''a18'', ''a38'', ''a14'', ''a12'', ''a40'', ''a24'', ''a25'', ''a43'', ''a44'', ''a2'', ''a10'', ''a8'', ''a20'', ''a47'', ''a37'', ''a35'', ''a50'', ''a39'', ''a46'', ''a45'', ''a28'', ''a15'', ''a19'']'
$wsRules.Range("B13").Value = '(age <=  32.0) => (class >= 3) [''a20'', ''a7'', ''a46'', ''a42'', ''a45'', ''a18'', ''a25'', ''a50'', ''a15'', ''a14'', ''a44'', ''a4'', ''a19'', ''a2'', ''a10'']'
$wsRules.Range("B15").Value = '(infertility <=  2.0) & (oocytes <=  5.0) => (class >= 3) [''a47'', ''a49'', ''a48'', ''a39'', ''a46'', ''a45'', ''a18'', ''a25'', ''a43'', ''a38'', ''a8'', ''a13'', ''a14'', ''a2'', ''a10'', ''a12'']'
$wsRules.Range("B16").Value = '(age <=  39.0) => (class >= 2) [''a26'', ''a9'', ''a18'', ''a33'', ''a38'', ''a34'', ''a14'', ''a29'', ''a12'', ''a30'', ''a40'', ''a17'', ''a24'', ''a25'', ''a32'', ''a43'', ''a21'', ''a4'', ''a44'', ''a2'', ''a10'', ''a8'', ''a20'', ''a47'', ''a37'', ''a36'', ''a16'', ''a31'', ''a35'', ''a50'', ''a7'', ''a11'', ''a39'', ''a46'', ''a42'', ''a45'', ''a28'', ''a41'', ''a1'', ''a27'', ''a3'', ''a15'', ''a13'', ''a19'', ''a22'']'
$wsRules.Range("B17").Value = '(woman_eval >=  3.0) => (class >= 2) [''a26'', ''a18'', ''a38'', ''a34'', ''a14'', ''a51'', ''a29'', ''a12'', ''a30'', ''a40'', ''a17'', ''a24'', ''a25'', ''a32'', ''a43'', ''a5'', ''a21'', ''a4'', ''a44'', ''a2'', ''a10'', ''a8'', ''a20'', ''a47'', ''a48'', ''a6'', ''a37'', ''a36'', ''a16'', ''a50'', ''a7'', ''a11'', ''a39'', ''a45'', ''a41'', ''a1'', ''a27'', ''a13'', ''a19'']'
$wsRules.Range("B18").Value = '(infertility <=  2.0) => (class >= 2) [''a18'', ''a38'', ''a14'', ''a49'', ''a12'', ''a25'', ''a43'', ''a2'', ''a10'', ''a8'', ''a47'', ''a48'', ''a7'', ''a39'', ''a46'', ''a45'', ''a1'', ''a3'', ''a13'', ''a22'']'

# --- Sheet "Statystyki reguł": update coverage values (column C) ---
$wsStats = $wb.Worksheets.Item("Statystyki reguł")
$wsStats.Range("C5").Value = 0.7441860465116279
$wsStats.Range("C9").Value = 0.09302325581395349
$wsStats.Range("C10").Value = 0.6923076923076923
$wsStats.Range("C14").Value = 0.05128205128205128
$wsStats.Range("C17").Value = 0.78

# --- Sheet "Walidacja krzyżowa": relabel + update accuracy/correct/f1_score rows ---
$wsCv = $wb.Worksheets.Item("Walidacja krzyżowa")
$wsCv.Range("A1").Value = "accuracy"
$wsCv.Range("B1").Value = 0.7450980392156863
$wsCv.Range("A2").Value = "not_classified"
$wsCv.Range("B2").Value = 0
$wsCv.Range("A3").Value = "correct"
$wsCv.Range("B3").Value = 0.7450980392156863
$wsCv.Range("A4").Value = "f1_score"
$wsCv.Range("B4").Value = 0
